$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: fill in the missing Profitable / SellPrice / Price Change % values
# and flip Holding to false (trade closed).
$ws.Range("B3").Value = $false
$ws.Range("E3").Value = 104.43
$ws.Range("F3").Value = -0.39106830946463644
$ws.Range("G3").Value = $false

# New row 4: next trade's starting Principle only.
$ws.Range("C4").Value = 9998.23

# Column F widened slightly to fit the repeater's new values.
# (Excel's ColumnWidth char-units get re-quantized on write; 11.666666666666666
# is the input that round-trips to a stored OOXML width of exactly 12.5.)
$ws.Columns.Item(6).ColumnWidth = 11.666666666666666
